$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "873"
$ws.Range("C2").Style = "Normal"
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "1953687.79"
$ws.Range("D2").Style = "Normal"

$ws.Range("C4").NumberFormat = "@"
$ws.Range("C4").Value = "1035"
$ws.Range("C4").Style = "Normal"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "3687879.47"
$ws.Range("D4").Style = "Normal"

$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "674"
$ws.Range("C6").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2202707.78"
$ws.Range("D6").Style = "Normal"

$ws.Range("C16").NumberFormat = "@"
$ws.Range("C16").Value = "500"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1858274.13"
$ws.Range("D16").Style = "Normal"

$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = "338"
$ws.Range("C21").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1203741.00"
$ws.Range("D21").Style = "Normal"

$ws.Range("C29").NumberFormat = "@"
$ws.Range("C29").Value = "300"
$ws.Range("C29").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "774236.89"
$ws.Range("D29").Style = "Normal"

$ws.Range("C31").NumberFormat = "@"
$ws.Range("C31").Value = "583"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2459419.89"
$ws.Range("D31").Style = "Normal"

$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = "403"
$ws.Range("C33").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1420480.57"
$ws.Range("D33").Style = "Normal"

$ws.Range("C35").NumberFormat = "@"
$ws.Range("C35").Value = "17"
$ws.Range("C35").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "49932.00"
$ws.Range("D35").Style = "Normal"

$ws.Range("C36").NumberFormat = "@"
$ws.Range("C36").Value = "368"
$ws.Range("C36").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "883760.71"
$ws.Range("D36").Style = "Normal"

$ws.Range("C37").NumberFormat = "@"
$ws.Range("C37").Value = "217"
$ws.Range("C37").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "672577.04"
$ws.Range("D37").Style = "Normal"

$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "208"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "563671.14"
$ws.Range("D38").Style = "Normal"

$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "6"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "13500.00"
$ws.Range("D39").Style = "Normal"

$ws.Range("C40").NumberFormat = "@"
$ws.Range("C40").Value = "18"
$ws.Range("C40").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "45492.85"
$ws.Range("D40").Style = "Normal"

$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "197"
$ws.Range("C41").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "548555.23"
$ws.Range("D41").Style = "Normal"

$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "94"
$ws.Range("C42").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "438699.98"
$ws.Range("D42").Style = "Normal"

$ws.Range("C43").NumberFormat = "@"
$ws.Range("C43").Value = "140"
$ws.Range("C43").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "572972.25"
$ws.Range("D43").Style = "Normal"

$ws.Range("C67").NumberFormat = "@"
$ws.Range("C67").Value = "24"
$ws.Range("C67").Style = "Normal"
$ws.Range("D67").NumberFormat = "@"
$ws.Range("D67").Value = "98660.00"
$ws.Range("D67").Style = "Normal"

$ws.Range("C68").NumberFormat = "@"
$ws.Range("C68").Value = "19"
$ws.Range("C68").Style = "Normal"
$ws.Range("D68").NumberFormat = "@"
$ws.Range("D68").Value = "70496.05"
$ws.Range("D68").Style = "Normal"

$ws.Range("C69").NumberFormat = "@"
$ws.Range("C69").Value = "6"
$ws.Range("C69").Style = "Normal"
$ws.Range("D69").NumberFormat = "@"
$ws.Range("D69").Value = "37700.00"
$ws.Range("D69").Style = "Normal"

$ws.Range("C75").NumberFormat = "@"
$ws.Range("C75").Value = "392"
$ws.Range("C75").Style = "Normal"
$ws.Range("D75").NumberFormat = "@"
$ws.Range("D75").Value = "993906.70"
$ws.Range("D75").Style = "Normal"

$ws.Range("C77").NumberFormat = "@"
$ws.Range("C77").Value = "941"
$ws.Range("C77").Style = "Normal"
$ws.Range("D77").NumberFormat = "@"
$ws.Range("D77").Value = "3309240.26"
$ws.Range("D77").Style = "Normal"

$ws.Range("C78").NumberFormat = "@"
$ws.Range("C78").Value = "528"
$ws.Range("C78").Style = "Normal"
$ws.Range("D78").NumberFormat = "@"
$ws.Range("D78").Value = "1739225.47"
$ws.Range("D78").Style = "Normal"

